# Daily attendance processing - 2025-10-29 17:21:28
#
# For every row in the "Recorded By" column (G), when the cell text is of
# the form "<name>, System" (i.e. "System" is the second of exactly two
# comma-separated entries), swap the order to "System, <name>".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string]) {
        $parts = $val -split ', '
        if ($parts.Count -eq 2 -and $parts[1] -eq 'System') {
            $cell.Value2 = "System, " + $parts[0]
        }
    }
}
